$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.844.41'
$ws.Range('E2').Value = '  -0.64%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.522.30'
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.32'
$ws.Range('E5').Value = '  -1.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '178.10'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.604'
$ws.Range('E8').Value = '  +0.57%  '
$ws.Range('B9').Value = 'LidoStakedEther'
$ws.Range('C9').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.516.14'
$ws.Range('E9').Value = '  +0.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.135'
$ws.Range('E10').Value = '  -1.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.93'
$ws.Range('E11').Value = '  -1.93%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.425'
$ws.Range('E12').Value = '  -2.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.126.93'
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '30.48'
$ws.Range('E14').Value = '  -4.41%  '
$ws.Range('E15').Value = '  -2.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.916.20'
$ws.Range('E16').Value = '  -0.58%  '
$ws.Range('E17').Value = '  -1.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.524.39'
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.09'
$ws.Range('E19').Value = '  -2.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.04'
$ws.Range('E20').Value = '  -1.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '383.88'
$ws.Range('E21').Value = '  -1.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.86'
$ws.Range('E22').Value = '  -1.66%  '
$ws.Range('E23').Value = '  +1.19%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.37'
$ws.Range('E24').Value = '  -2.00%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.997'
$ws.Range('E25').Value = '  -0.45%  '
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.75'
$ws.Range('E26').Value = '  +0.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000121'
$ws.Range('E27').Value = '  +0.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.95'
$ws.Range('E28').Value = '  -3.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.173'
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.994'
$ws.Range('E30').Value = '  -0.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '24.53'
$ws.Range('E31').Value = '  +4.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.93'
$ws.Range('E32').Value = '  -3.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.02'
$ws.Range('E33').Value = '  -1.51%  '
$ws.Range('E34').Value = '  -3.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.27'
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '160.90'
$ws.Range('E38').Value = '  -2.26%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '29.32'
$ws.Range('E39').Value = '  +12.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.894'
$ws.Range('E40').Value = '  +2.95%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.80'
$ws.Range('E41').Value = '  -3.33%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.54'
$ws.Range('E42').Value = '  -2.39%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.57'
$ws.Range('E43').Value = '  -3.61%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.59'
$ws.Range('E44').Value = '  -5.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.745.65'
$ws.Range('E45').Value = '  -3.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0707'
$ws.Range('E46').Value = '  -2.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '25.47'
$ws.Range('E47').Value = '  -5.74%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '40.73'
$ws.Range('E48').Value = '  -2.08%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0298'
$ws.Range('E49').Value = '  -0.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '324.99'
$ws.Range('E50').Value = '  -2.90%  '
$ws.Range('E51').Value = '  -2.47%  '
